# Fix numba error in Output_Centre.py
#
# The "Values" row (row 4) of the Numerical Method column used "ER" as the
# numerical method, which triggers a numba error in Output_Centre.py.
# Switch it to "V" instead (the other accepted value per the "Accepted
# values" row: "string (ER or V)").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("M4").Value = "V"

# Reflect the view/selection state recorded alongside the edit: the user's
# selection ended up on M8 with the window scrolled so column F is the
# left-most visible column.
$win = $excel.ActiveWindow
$ws.Range("M8").Select()
try { $win.ScrollColumn = 6 } catch {}
try { $win.ScrollRow = 1 } catch {}
